$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.787.73"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "2.091.61"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("E6").Value = "  -1.52%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("E10").Value = "  -3.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.909"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.13%  "

$ws.Range("D15").Value = "2.397.02"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.48%  "

$ws.Range("D17").Value = "2.104.64"
$ws.Range("E17").Value = "  +2.51%  "

$ws.Range("D18").Value = "36.749.97"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").Value = "0.0₃0881"
$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  -3.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.22%  "

$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.88%  "

$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("E35").Value = "  +7.46%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.88%  "

$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0953"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.88%  "

$ws.Range("D47").Value = "1.389.18"
$ws.Range("E47").Value = "  +9.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").Value = "2.283.87"
$ws.Range("E51").Value = "  +1.92%  "
